$d = $word.ActiveDocument

# 1) " GARIS BESAR HALUAN ORGANISASI" + " " (two runs) -> " GARIS BESAR HALUAN ORGANISASI " (one run)
#    The search string spans both existing runs (identical bold/lang formatting), so Word
#    rewrites the matched range as a single run carrying the combined text.
$d.Content.Find.Execute(" GARIS BESAR HALUAN ORGANISASI ", $true, $false, $false, $false, $false, `
    $true, 1, $false, " GARIS BESAR HALUAN ORGANISASI ", 2) | Out-Null

# 2) "Pukul : 23.03 WIB" -> "Pukul : 23.14 WIB" -- only the minutes run ("03") changes.
#    "03" is unique across the whole document body, so this only touches that single run.
$d.Content.Find.Execute("03", $true, $false, $false, $false, $false, `
    $true, 1, $false, "14", 2) | Out-Null

Write-Output "done"
